$d = $word.ActiveDocument

function InsertBreakBefore($searchText) {
    $r = $d.Content
    $r.Find.Execute($searchText, $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0) | Out-Null
    $r.Collapse(1)  # wdCollapseStart
    $r.InsertParagraphBefore()
}

function InsertBreakAfter($searchText) {
    $r = $d.Content
    $r.Find.Execute($searchText, $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0) | Out-Null
    $r.Collapse(0)  # wdCollapseEnd
    $r.InsertParagraphAfter()
}

# Break the single big paragraph into one paragraph per run (9 new paragraph
# boundaries across the 10 original runs). Use "^l" in the find text to make
# sure manual line breaks (<w:br/>) stay attached to the run/paragraph that
# precedes them, matching the original run grouping.
InsertBreakBefore("Fig.3 Scientific Articles from year 2009")
InsertBreakAfter("entific Articles from year 2009 to 2021 ")
InsertBreakAfter(" the concept of homomorphic encryption. ^l")
InsertBreakAfter("Inclusion and exclusion criteria ^l")
InsertBreakAfter("hcare or bioinfor matics were excluded. ")
InsertBreakAfter("Research questions ^l")
InsertBreakAfter("s systematic literature review in Table ")
InsertBreakBefore("Background ")
InsertBreakAfter("Background ^l")

# Switch page size from Letter to A4 (keeps the 2-column layout/margins).
$ps = $d.PageSetup
$ps.PageWidth = 11906 / 20.0
$ps.PageHeight = 16838 / 20.0
